$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header / shared-string text values ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Update data values (column C: GDP) ---
$ws.Range("C2").Value = 5191.140356354663
$ws.Range("C3").Value = 5082.354756663512
$ws.Range("C4").Value = 5660.517066940175
$ws.Range("C5").Value = 5360.226632400601
$ws.Range("C6").Value = 5642.578115155247
$ws.Range("C7").Value = 5919.20956823756
$ws.Range("C8").Value = 5996.49696468919
$ws.Range("C9").Value = 6301.696269820412
$ws.Range("C10").Value = 6114.227214287786
$ws.Range("C11").Value = 6661.86504232374
$ws.Range("C12").Value = 6262.368904654469

# --- Update data values (column AL: Delegation flag) ---
$ws.Range("AL2").Value = 1
$ws.Range("AL4").Value = 1
$ws.Range("AL9").Value = 1
$ws.Range("AL11").Value = 1
